$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, $Text)
    $origStyle = $Cell.Style
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.Style = $origStyle
}

Set-TextValue $ws.Range('D2') '246.05'
Set-TextValue $ws.Range('D3') '22.00'
Set-TextValue $ws.Range('D4') '5.379'
Set-TextValue $ws.Range('D5') '0.05803'
Set-TextValue $ws.Range('D6') '3.377'
Set-TextValue $ws.Range('D8') '0.8062'
Set-TextValue $ws.Range('D9') '1.003'
Set-TextValue $ws.Range('B10') 'One'
Set-TextValue $ws.Range('C10') 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
Set-TextValue $ws.Range('D10') '0.01120'
Set-TextValue $ws.Range('E10') '9OneONEBestin24h'
Set-TextValue $ws.Range('B11') 'WazirX'
Set-TextValue $ws.Range('C11') 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextValue $ws.Range('D11') '0.1427'
Set-TextValue $ws.Range('E11') '10WazirXWRX'
Set-TextValue $ws.Range('B12') 'MandalaExchangeToken'
Set-TextValue $ws.Range('C12') 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextValue $ws.Range('D12') '0.07513'
Set-TextValue $ws.Range('E12') '11MandalaExchangeTokenMDX'
Set-TextValue $ws.Range('B13') 'LiechtensteinCryptoassetsExchange'
Set-TextValue $ws.Range('C13') 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextValue $ws.Range('D13') '0.03194'
Set-TextValue $ws.Range('E13') '12LiechtensteinCryptoassetsExchangeLCX'
Set-TextValue $ws.Range('B14') 'BitrueCoin'
Set-TextValue $ws.Range('C14') 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextValue $ws.Range('D14') '0.03057'
Set-TextValue $ws.Range('E14') '13BitrueCoinBTR'
Set-TextValue $ws.Range('B15') 'MCDex'
Set-TextValue $ws.Range('C15') 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
Set-TextValue $ws.Range('D15') '4.165'
Set-TextValue $ws.Range('E15') '14MCDexMCB'
Set-TextValue $ws.Range('B16') 'BitMartToken'
Set-TextValue $ws.Range('C16') 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextValue $ws.Range('D16') '0.09398'
Set-TextValue $ws.Range('E16') '15BitMartTokenBMX'
Set-TextValue $ws.Range('B17') 'BitForexToken'
Set-TextValue $ws.Range('C17') 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextValue $ws.Range('D17') '0.001601'
Set-TextValue $ws.Range('E17') '16BitForexTokenBF'
Set-TextValue $ws.Range('B18') 'CoinExToken'
Set-TextValue $ws.Range('C18') 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
Set-TextValue $ws.Range('D18') '0.04800'
Set-TextValue $ws.Range('E18') '17CoinExTokenCET'
Set-TextValue $ws.Range('D19') '0.005650'
Set-TextValue $ws.Range('D20') '0.004094'
Set-TextValue $ws.Range('E20') '19HotbitTokenHTBWorstin24h'
Set-TextValue $ws.Range('D21') '0.0009965'
Set-TextValue $ws.Range('D23') '3.699'
Set-TextValue $ws.Range('D24') '2.246'
Set-TextValue $ws.Range('D25') '0.3204'
Set-TextValue $ws.Range('D27') '0.0003592'
Set-TextValue $ws.Range('E27') '26UpBotsUBXT'
Set-TextValue $ws.Range('D40') '0.03887'
Set-TextValue $ws.Range('B41') 'KickToken'
Set-TextValue $ws.Range('C41') 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
Set-TextValue $ws.Range('D41') '0.006319'
Set-TextValue $ws.Range('E41') '40KickTokenKICK'
Set-TextValue $ws.Range('B42') 'BKEXToken'
Set-TextValue $ws.Range('C42') 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
Set-TextValue $ws.Range('D42') '0.1073'
Set-TextValue $ws.Range('E42') '41BKEXTokenBKK'
Set-TextValue $ws.Range('B43') 'CEJI'
Set-TextValue $ws.Range('C43') 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
Set-TextValue $ws.Range('D43') '0.003000'
Set-TextValue $ws.Range('E43') '42CEJICEJI'
Set-TextValue $ws.Range('D44') '0.006685'
Set-TextValue $ws.Range('D45') '0.00005591'
Set-TextValue $ws.Range('D47') '0.3899'
Set-TextValue $ws.Range('D48') '0.1448'
